$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  46008, 10819.5183127189, 9686.86193830566, 16512.26, 6851.81245724934, 1.10059981479162),
    @(3,  46009, 10620.3379674477, 9665.35760894511, 11232.26, 6949.08627034535, 224.257661637103),
    @(4,  46010, 10637.1481954349, 9054.20847595399, 11232.26, 6962.73184955101, 199.361680229375),
    @(5,  46011, 4127.66864886361, 6365.40484883612, 11232.26, 6707.48663969339, 76.6929786887295),
    @(6,  46012, 4136.85413087738, 6492.62014812714, 11232.26, 6890.38630530073, 89.6144355594946),
    @(7,  46013, 9566.18303866457, 9255.11389491146, 11232.26, 8113.18340118371, 255.668220670632),
    @(8,  46014, 9566.18303866457, 9451.54200333216, 17376.26, 8113.18340118371, 7.85272518816131),
    @(9,  46015, 9566.18303866457, 9541.862540334,   17208.26, 8113.18340118371, 18.6160808965712),
    @(10, 46016, 8471.80799626025, 8538.06430671732, 16416.26, 7768.80692669759, -4.55786527437886),
    @(11, 46017, 8471.80799626025, 7886.16588242443, 15840.26, 7768.80692669759, -7.72029961991598),
    @(12, 46018, 8547.24893097507, 8590.47310213652, 16632.26, 7775.98586474824, -11.0750430464682),
    @(13, 46019, 8471.80799626025, 8695.06866928379, 17448.26, 7768.80692669759, -41.0160168341091),
    @(14, 46020, 9603.52690270833, 9419.78850845279, 19176.26, 8403.10977092558, -56.3900716925677),
    @(15, 46021, 9603.52690270833, 9583.7804313854,  11232.26, 8403.10977092558, 281.442925096291)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
